$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate Species / Disease / Symptoms columns (rows 2-25) from English to French
$ws.Range("A2").Value2 = 'volaille'
$ws.Range("B2").Value2 = 'La grippe aviaire'
$ws.Range("C2").Value2 = 'Fièvre, Mal de gorge, Toux, Maux de tête, Douleur musculaire, Respiration difficile, Conjonctivite (la_partie_intérieure_de_la_paupière)'

$ws.Range("A3").Value2 = 'bétail'
$ws.Range("B3").Value2 = 'fièvre de la vallée du Rift'
$ws.Range("C3").Value2 = 'Fièvre, léthargie, mort subite, écoulement nasal, salivation excessive, anorexie, diarrhée'

$ws.Range("A4").Value2 = 'bétail'
$ws.Range("B4").Value2 = 'Tuberculose bovine'
$ws.Range("C4").Value2 = 'état général de la maladie, perte de poids, mort subite'

$ws.Range("A5").Value2 = 'bétail'
$ws.Range("B5").Value2 = 'Peau grumeleuse'
$ws.Range("C5").Value2 = 'Nodules sur la peau, Peau lésée'

$ws.Range("A6").Value2 = 'bétail'
$ws.Range("B6").Value2 = 'Peste des Petits Ruminants'
$ws.Range("C6").Value2 = 'Fièvre, lésions buccales, diarrhée, mort subite'

$ws.Range("A7").Value2 = 'bétail'
$ws.Range("B7").Value2 = 'Diarrhée virale bovine'
$ws.Range("C7").Value2 = 'Fièvre, léthargie, anorexie, sécrétions oculaires, écoulement nasal, lésions buccales, diarrhée, réduction de la production'

$ws.Range("A8").Value2 = 'volaille'
$ws.Range("B8").Value2 = 'Newcastle'
$ws.Range("C8").Value2 = 'Sneezing, Nasal discharge, Cough, Diarrhea, Shudder, Drooping wings, Paralysis, Swelling of the tissues around the eyes and the neck, sudden death,  Reduction in production'

$ws.Range("A9").Value2 = 'volaille'
$ws.Range("B9").Value2 = 'Coryza infectieux'
$ws.Range("C9").Value2 = 'Tête ou visage enflé, éternuements, toux, sécrétions oculaires, écoulement nasal, anorexie, respiration difficile'

$ws.Range("A10").Value2 = 'volaille'
$ws.Range("B10").Value2 = 'Muguet'
$ws.Range("C10").Value2 = 'Anorexie, lésions buccales, croissance lente, yeux squameux, ailes tombantes, respiration difficile, respiration bruyante'

$ws.Range("A11").Value2 = 'bétail'
$ws.Range("B11").Value2 = 'Mastite (infection bactérienne)'
$ws.Range("C11").Value2 = 'Taille anormale, dureté de la mamelle, fièvre, changement de couleur du lait'

$ws.Range("A12").Value2 = 'bétail'
$ws.Range("B12").Value2 = 'Tuberculose'
$ws.Range("C12").Value2 = 'Abcès, toux, ganglions lymphatiques enflés, augmentation de la fréquence cardiaque'

$ws.Range("A13").Value2 = 'bétail'
$ws.Range("B13").Value2 = 'Gonfler'
$ws.Range("C13").Value2 = 'Flatulences, salivation excessive, gémissements, anorexie, vomissements'

$ws.Range("A14").Value2 = 'bétail'
$ws.Range("B14").Value2 = 'Encéphalomyélite équine'
$ws.Range("C14").Value2 = 'Fièvre, vision double, démarche irrégulière, frissons, lésions buccales, marche sans but'

$ws.Range("A15").Value2 = 'bétail'
$ws.Range("B15").Value2 = 'Septicémie hémorragique'
$ws.Range("C15").Value2 = 'Fièvre, respiration difficile, toux, sécrétions oculaires, écoulement nasal'

$ws.Range("A16").Value2 = 'bétail'
$ws.Range("B16").Value2 = 'Brucellose'
$ws.Range("C16").Value2 = 'Gonflement des testicules, bactéries localisées dans les articulations, Fièvre, Frissons, Anorexie, Transpiration, Léthargie'

$ws.Range("A17").Value2 = 'bétail'
$ws.Range("B17").Value2 = 'Clavelée'
$ws.Range("C17").Value2 = 'Fièvre, peau endommagée, inflammation de la muqueuse nasale, lésions cutanées, sécrétions oculaires, écoulement nasal, paupières enflées, léthargie, anorexie, paralysie'

$ws.Range("A18").Value2 = 'bétail'
$ws.Range("B18").Value2 = 'Pneumonie'
$ws.Range("C18").Value2 = 'Fièvre, anorexie, respiration difficile, sécrétions oculaires, écoulement nasal, salivation excessive, diarrhée'

$ws.Range("A19").Value2 = 'bétail'
$ws.Range("B19").Value2 = 'Fièvre aphteuse (FA)'
$ws.Range("C19").Value2 = 'Fièvre, lésions cutanées, peau endommagée, respiration difficile, salivation excessive, lésions buccales'

$ws.Range("A20").Value2 = 'bétail'
$ws.Range("B20").Value2 = 'Fièvre catarrhale du mouton'
$ws.Range("C20").Value2 = 'Fièvre, Lésions buccales, Respiration difficile, Langue violacée, Boiterie'

$ws.Range("A21").Value2 = 'bétail'
$ws.Range("B21").Value2 = 'Anaplasmose'
$ws.Range("C21").Value2 = 'Fièvre, pâleur autour des yeux, léthargie, perte de poids, réduction de la production, comportement agressif'

$ws.Range("A22").Value2 = 'bétail'
$ws.Range("B22").Value2 = 'Rage'
$ws.Range("C22").Value2 = 'Anorexie, prurit, boiterie, ténesme, salivation excessive, comportement agressif'

$ws.Range("A23").Value2 = 'volaille'
$ws.Range("B23").Value2 = 'Entérite nécrotique'
$ws.Range("C23").Value2 = 'Anorexie, Léthargie, Plumes pelucheuses, Yeux fermés, Diarrhée'

$ws.Range("A24").Value2 = 'volaille'
$ws.Range("B24").Value2 = 'Ascaris'
$ws.Range("C24").Value2 = 'Anorexie, Diarrhée, Croissance lente, Léthargie, Plumes ébouriffées, Perte de poids, changements de comportement'

$ws.Range("A25").Value2 = 'volaille'
$ws.Range("B25").Value2 = 'Variole aviaire'
$ws.Range("C25").Value2 = 'Paupières enflées, Yeux fermés, Lésions buccales, Perte de poids, Anorexie'

# Set explicit column widths
$ws.Columns.Item(1).ColumnWidth = 20.8
$ws.Columns.Item(2).ColumnWidth = 30
$ws.Columns.Item(3).ColumnWidth = 136.5

# Move the active selection to C25
$ws.Range("C25").Select()

